$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 3624.75
$ws.Range("I28").Value = 3624.75
$ws.Range("K28").Value = 3624.75
$ws.Range("M28").Value = -3139.75
# Row 98
$ws.Range("H98").Value = 697.6
$ws.Range("I98").Value = 697.6
$ws.Range("K98").Value = 697.6
$ws.Range("M98").Value = 800.4
# Row 122
$ws.Range("H122").Value = 697.6
$ws.Range("I122").Value = 697.6
$ws.Range("K122").Value = 2092.8
$ws.Range("M122").Value = 357.1999999999998
# Row 125
$ws.Range("H125").Value = 4361.5713
$ws.Range("I125").Value = 2766
$ws.Range("J125").Value = 4999.8
$ws.Range("K125").Value = 24894
$ws.Range("L125").Value = 44998.2
$ws.Range("M125").Value = -22434
$ws.Range("N125").Value = -49918.2
# Row 138
$ws.Range("H138").Value = 2712.6072
$ws.Range("I138").Value = 1124.1111
$ws.Range("K138").Value = 3372.3333
$ws.Range("M138").Value = 1767.6667

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
# Row 45
$ws.Range("H45").Value = 1949.8
$ws.Range("I45").Value = 1562.25
$ws.Range("K45").Value = 1562.25
$ws.Range("M45").Value = -1185.25
# Row 61
$ws.Range("H61").Value = 2491.4375
$ws.Range("I61").Value = 2457.7334
$ws.Range("K61").Value = 2457.7334
$ws.Range("M61").Value = -2245.7334
# Row 97
$ws.Range("H97").Value = 809.8125
$ws.Range("I97").Value = 401.07693
$ws.Range("K97").Value = 401.07693
$ws.Range("M97").Value = 94.92307
# Row 116
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = ""
$ws.Range("N116").Value = ""
# Row 122
$ws.Range("H122").Value = 1492
$ws.Range("I122").Value = 1492
$ws.Range("K122").Value = 4476
$ws.Range("M122").Value = -2026
# Row 132
$ws.Range("H132").Value = 2603.111
$ws.Range("I132").Value = 2186.3914
$ws.Range("J132").Value = 4999.25
$ws.Range("K132").Value = 6559.174199999999
$ws.Range("L132").Value = 14997.75
$ws.Range("M132").Value = -4029.174199999999
$ws.Range("N132").Value = -20057.75
# Row 136
$ws.Range("H136").Value = 2491.4375
$ws.Range("I136").Value = 2457.7334
$ws.Range("K136").Value = 7373.2002
$ws.Range("M136").Value = -4823.2002
# Row 141
$ws.Range("H141").Value = 89999
$ws.Range("J141").Value = 89999
$ws.Range("L141").Value = 89999
$ws.Range("N141").Value = -100359

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
# Row 20
$ws.Range("H20").Value = 2931.1765
$ws.Range("I20").Value = 2903.8333
$ws.Range("K20").Value = 2903.8333
$ws.Range("M20").Value = -2656.8333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1850.1428
$ws.Range("I16").Value = 1790.4
$ws.Range("K16").Value = 1790.4
$ws.Range("M16").Value = -1503.4
# Row 31
$ws.Range("H31").Value = 4690.6
$ws.Range("I31").Value = 2705.7144
$ws.Range("J31").Value = 9322
$ws.Range("K31").Value = 2705.7144
$ws.Range("L31").Value = 9322
$ws.Range("M31").Value = -2410.7144
$ws.Range("N31").Value = -9912
# Row 34
$ws.Range("H34").Value = 4690.6
$ws.Range("I34").Value = 2705.7144
$ws.Range("J34").Value = 9322
$ws.Range("K34").Value = 2705.7144
$ws.Range("L34").Value = 9322
$ws.Range("M34").Value = -2503.7144
$ws.Range("N34").Value = -9726
# Row 113
$ws.Range("H113").Value = 1850.1428
$ws.Range("I113").Value = 1790.4
$ws.Range("K113").Value = 1790.4
$ws.Range("M113").Value = 379.5999999999999
# Row 122
$ws.Range("H122").Value = 1482.6875
$ws.Range("I122").Value = 1652
$ws.Range("J122").Value = 974.75
$ws.Range("K122").Value = 4956
$ws.Range("L122").Value = 2924.25
$ws.Range("M122").Value = -2506
$ws.Range("N122").Value = -7824.25
# Row 132
$ws.Range("H132").Value = 2905.2334
$ws.Range("I132").Value = 2683.0386
$ws.Range("J132").Value = 4349.5
$ws.Range("K132").Value = 8049.1158
$ws.Range("L132").Value = 13048.5
$ws.Range("M132").Value = -5519.1158
$ws.Range("N132").Value = -18108.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 774.5
$ws.Range("I97").Value = 749
$ws.Range("K97").Value = 2247
$ws.Range("M97").Value = -1751
# Row 114
$ws.Range("H114").Value = 6604.4
$ws.Range("I114").Value = 4341.6665
$ws.Range("J114").Value = 9998.5
$ws.Range("K114").Value = 13024.9995
$ws.Range("L114").Value = 29995.5
$ws.Range("M114").Value = -9770.999500000002
$ws.Range("N114").Value = -36503.5
# Row 117
$ws.Range("H117").Value = 1140.3334
$ws.Range("I117").Value = 489
$ws.Range("J117").Value = 1466
$ws.Range("K117").Value = 1467
$ws.Range("L117").Value = 4398
$ws.Range("M117").Value = 1975
$ws.Range("N117").Value = -11282
# Row 140
$ws.Range("H140").Value = 3162.5
$ws.Range("I140").Value = 3095.2
$ws.Range("K140").Value = 9285.599999999999
$ws.Range("M140").Value = -4105.599999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 774
$ws.Range("I122").Value = 774
$ws.Range("J122").Value = 774
$ws.Range("K122").Value = 2322
$ws.Range("L122").Value = 2322
$ws.Range("M122").Value = 128
$ws.Range("N122").Value = -7222
# Row 123
$ws.Range("H123").Value = 48713
$ws.Range("J123").Value = 48713
$ws.Range("L123").Value = 48713
$ws.Range("N123").Value = -53613
# Row 132
$ws.Range("H132").Value = 3328.75
$ws.Range("I132").Value = 2993.7144
$ws.Range("J132").Value = 3797.8
$ws.Range("K132").Value = 8981.143199999999
$ws.Range("L132").Value = 11393.4
$ws.Range("M132").Value = -6451.143199999999
$ws.Range("N132").Value = -16453.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2802.6667
$ws.Range("I7").Value = 2802.6667
$ws.Range("K7").Value = 2802.6667
$ws.Range("M7").Value = -2690.6667
# Row 40
$ws.Range("H40").Value = 5163.3335
$ws.Range("I40").Value = 5163.3335
$ws.Range("K40").Value = 5163.3335
$ws.Range("M40").Value = -5027.3335
# Row 126
$ws.Range("H126").Value = 2802.6667
$ws.Range("I126").Value = 2802.6667
$ws.Range("K126").Value = 8408.000100000001
$ws.Range("M126").Value = -5938.000100000001
# Row 136
$ws.Range("H136").Value = 2876.5293
$ws.Range("I136").Value = 2962.5625
$ws.Range("K136").Value = 8887.6875
$ws.Range("M136").Value = -6337.6875
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 24000
$ws.Range("I62").Value = 24000
$ws.Range("K62").Value = 24000
$ws.Range("M62").Value = -23376
# Row 65
$ws.Range("H65").Value = 24000
$ws.Range("I65").Value = 24000
$ws.Range("K65").Value = 120000
$ws.Range("M65").Value = -116880
# Row 100
$ws.Range("H100").Value = 1435.75
$ws.Range("I100").Value = 2300
$ws.Range("J100").Value = 1147.6666
$ws.Range("K100").Value = 4600
$ws.Range("L100").Value = 2295.3332
$ws.Range("M100").Value = -4059
$ws.Range("N100").Value = -3377.3332
# Row 132
$ws.Range("H132").Value = 1728.8182
$ws.Range("I132").Value = 1263.1666
$ws.Range("J132").Value = 3824.25
$ws.Range("K132").Value = 3789.4998
$ws.Range("L132").Value = 11472.75
$ws.Range("M132").Value = -1259.4998
$ws.Range("N132").Value = -16532.75
# Row 136
$ws.Range("H136").Value = 712.1539
$ws.Range("I136").Value = 712.1539
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2136.4617
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 413.5383000000002
$ws.Range("N136").Value = ""
